$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.595.47"
$ws.Range("E2").Value = "  +0.84%  "

$ws.Range("D3").Value = "1.639.43"
$ws.Range("E3").Value = "  +1.08%  "

$ws.Range("E4").Value = "  -0.29%  "

$ws.Range("D5").Value = "'214.20"
$ws.Range("E5").Value = "  +0.92%  "

$ws.Range("E6").Value = "  +1.66%  "

$ws.Range("E7").Value = "  -0.22%  "

$ws.Range("E8").Value = "  +1.16%  "

$ws.Range("D9").Value = "'0.0624"
$ws.Range("E9").Value = "  +0.76%  "

$ws.Range("E10").Value = "  +1.09%  "

$ws.Range("D11").Value = "'0.0840"
$ws.Range("E11").Value = "  -0.14%  "

$ws.Range("D12").Value = "1.867.77"
$ws.Range("E12").Value = "  +0.95%  "

$ws.Range("D13").Value = "1.645.14"
$ws.Range("E13").Value = "  +1.20%  "

$ws.Range("E14").Value = "  +1.78%  "

$ws.Range("E15").Value = "  +1.49%  "

$ws.Range("D16").Value = "'64.77"
$ws.Range("E16").Value = "  +1.14%  "

$ws.Range("D17").Value = "26.598.84"
$ws.Range("E17").Value = "  +0.79%  "

$ws.Range("E18").Value = "  +0.49%  "

$ws.Range("D19").Value = "'214.96"
$ws.Range("E19").Value = "  +0.38%  "

$ws.Range("E20").Value = "  -0.13%  "

$ws.Range("E21").Value = "  +0.87%  "

$ws.Range("D22").Value = "'6.23"
$ws.Range("E22").Value = "  +0.67%  "

$ws.Range("D23").Value = "'9.43"
$ws.Range("E23").Value = "  +1.91%  "

$ws.Range("E24").Value = "  +12.96%  "

$ws.Range("D25").Value = "'144.78"
$ws.Range("E25").Value = "  -1.90%  "

$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("E27").Value = "  +0.05%  "

$ws.Range("E28").Value = "  +4.38%  "

$ws.Range("E29").Value = "  +0.82%  "

$ws.Range("E30").Value = "  +1.00%  "

$ws.Range("E31").Value = "  +1.28%  "

$ws.Range("E32").Value = "  +1.26%  "

$ws.Range("D34").Value = "1.274.66"
$ws.Range("E34").Value = "  +5.22%  "

$ws.Range("E35").Value = "  +2.84%  "

$ws.Range("D36").Value = "'2.40"
$ws.Range("E36").Value = "  +1.18%  "

$ws.Range("E37").Value = "  +2.63%  "

$ws.Range("E38").Value = "  +6.37%  "

$ws.Range("E39").Value = "  +3.87%  "

$ws.Range("E40").Value = "  -0.16%  "

$ws.Range("D41").Value = "'0.808"
$ws.Range("E41").Value = "  +2.11%  "

$ws.Range("E42").Value = "  -0.07%  "

$ws.Range("D43").Value = "'5.40"
$ws.Range("E43").Value = "  +0.79%  "

$ws.Range("D44").Value = "1.778.55"
$ws.Range("E44").Value = "  +1.10%  "

$ws.Range("D45").Value = "'91.13"
$ws.Range("E45").Value = "  -1.50%  "

$ws.Range("D46").Value = "'59.31"
$ws.Range("E46").Value = "  +8.67%  "

$ws.Range("E47").Value = "  +1.27%  "

$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0101"
$ws.Range("E48").Value = "  -0.66%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.0515"
$ws.Range("E49").Value = "  +0.92%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'7.69"
$ws.Range("E50").Value = "  +1.23%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.0960"
$ws.Range("E51").Value = "  +1.27%  "

Write-Output "Updated cryptos list"